# CrudOperation.xlsx / Address Book "File path" test-data refresh
# Sheet1 (Address Book) columns P/Q/R hold the per-row FedEx tracking
# number (P), expected/actual rate (Q) and pass/fail result (R) for the
# smoke-test rows. This updates the stale tracking numbers / rate /
# result to the new values recorded against the refreshed shipment file.
#
# The numeric-looking tracking numbers and the "$"-prefixed rates must
# stay stored as plain text (as they were originally), so they are
# entered with a leading apostrophe (quote-prefix) just like typing them
# into Excel would do, and the style is then reset back to "Normal" so
# no new number-format/style gets left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: tracking number only changes
$ws.Range("P2").Value = "'320018612583"
$ws.Range("P2").Style = "Normal"

# Row 3: tracking number only changes
$ws.Range("P3").Value = "'320018612594"
$ws.Range("P3").Style = "Normal"

# Row 4: tracking number, actual rate, and result all change (FAIL -> PASS)
$ws.Range("P4").Value = "'320018612620"
$ws.Range("P4").Style = "Normal"

$ws.Range("Q4").Value = "'`$49.70"
$ws.Range("Q4").Style = "Normal"

$ws.Range("R4").Value = "PASS"

# Row 5: tracking number and actual rate change
$ws.Range("P5").Value = "'320018607974"
$ws.Range("P5").Style = "Normal"

$ws.Range("Q5").Value = "'`$43.36"
$ws.Range("Q5").Style = "Normal"
